$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new 2021 column (R), mirroring the formatting of
# the existing 2020 column (Q) for each row.
$ws.Range("Q4:Q14").Copy($ws.Range("R4:R14"))

# Overwrite the copied values with the actual 2021 figures.
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 6.5159856023713738
$ws.Range("R6").Value = 25.411968777103212
$ws.Range("R7").Value = 4.5359966708281316
$ws.Range("R8").Value = 9.213483146067416
$ws.Range("R9").Value = 12.204234122042342
$ws.Range("R10").Value = 9.4037615046018406
$ws.Range("R11").Value = 5.6537102473498235
$ws.Range("R12").Value = 1.5984015984015985
$ws.Range("R13").Value = 6.2881802387490886
$ws.Range("R14").Value = 8.1261101243339251

# Move / update the active selection to match the post-edit state.
$ws.Range("S8").Select()
